$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.446.44'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.379.83'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.78'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.38'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.99%  '
$ws.Range('E6').ClearFormats()

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.636'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('E7').ClearFormats()

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.617'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.37%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.18'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.53%  '
$ws.Range('E10').ClearFormats()

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0921'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('E11').ClearFormats()

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.56'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('E12').ClearFormats()

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.110'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('E13').ClearFormats()

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.39%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.742.32'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.53'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.378.86'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.388.35'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('E18').ClearFormats()

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.66'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +15.44%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.35'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.69%  '
$ws.Range('E20').ClearFormats()

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('E21').ClearFormats()

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.59%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.36'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.93%  '
$ws.Range('E23').ClearFormats()

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '261.59'
$ws.Range('D24').ClearFormats()

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.97%  '
$ws.Range('E24').ClearFormats()

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E25').ClearFormats()

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.59'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.30'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E28').ClearFormats()

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.62%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.46'
$ws.Range('D30').ClearFormats()

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0955'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E31').ClearFormats()

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '37.60'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.39%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '167.88'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.92%  '
$ws.Range('E33').ClearFormats()

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.91'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.42%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.132'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E35').ClearFormats()

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('E36').ClearFormats()

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.76'
$ws.Range('D37').ClearFormats()

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.57%  '
$ws.Range('E37').ClearFormats()

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.97'
$ws.Range('D38').ClearFormats()

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +14.65%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.04'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.86%  '
$ws.Range('E39').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.98'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('E40').ClearFormats()

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0356'
$ws.Range('D41').ClearFormats()

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.72'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -7.15%  '
$ws.Range('E42').ClearFormats()

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.53'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.29%  '
$ws.Range('E43').ClearFormats()

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.97%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.00'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.09%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.855.28'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +12.97%  '
$ws.Range('E46').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('E47').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.31%  '
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '83.72'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +6.26%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '113.01'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.15%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.29'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.80%  '
$ws.Range('E51').ClearFormats()

